# This script applies targeted text replacements while preserving the
# exact run structure of each paragraph (including any pre-existing empty
# <w:r/> runs), which a plain Find.Execute-based replace would otherwise
# collapse/merge away in this runtime. It does this by using
# Range.InsertXML to swap the *existing* text-bearing run's XML in place
# (rather than deleting the old text and then inserting new text, which
# triggers the run-list to be rebuilt and any zero-length sibling runs to
# be dropped).

function Replace-ParaRun($d, $startPos, $endPos, $newText, $rPr) {
    # Replace the run spanning [$startPos, $endPos) -- i.e. everything in
    # the paragraph up to, but excluding, the paragraph mark -- with a
    # single new run carrying $newText and the (optional) run properties
    # XML fragment $rPr. Any other runs already in the paragraph (such as a
    # leading/trailing empty <w:r/>) are left completely untouched.
    $rng = $d.Range($startPos, $endPos - 1)
    $runXml = "<w:r>" + $rPr + "<w:t>" + $newText + "</w:t></w:r>"
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

$d = $word.ActiveDocument

# Snapshot every paragraph's position/text/style up front so that later
# edits (which shift character offsets) cannot disturb earlier lookups.
$paras = @()
foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    $styleName = ""
    try { $styleName = $p.Style.NameLocal } catch {}
    $paras += [PSCustomObject]@{
        Start = $rng.Start
        End   = $rng.End
        Text  = $rng.Text.TrimEnd()
        Style = $styleName
    }
}

# Each entry: old paragraph text, new paragraph text, and an optional
# paragraph style name used to disambiguate cases where the same old text
# occurs more than once (here, the title appears both as the Heading 1 and
# as a bold "Normal" paragraph near the end, and they must keep their own
# formatting).
$edits = @(
    [PSCustomObject]@{ Old = "Play Aloha! Cluster Pays Free: Review and Strategies"; New = "Play Aloha! Cluster Pays for Free - Exciting Online Slot Game"; Style = "Heading 1"; RPr = "" }
    [PSCustomObject]@{ Old = "Play Aloha! Cluster Pays Free: Review and Strategies"; New = "Play Aloha! Cluster Pays for Free - Exciting Online Slot Game"; Style = "Normal"; RPr = "<w:rPr><w:b/></w:rPr>" }
    [PSCustomObject]@{ Old = "Unique cluster pays mechanism"; New = "Unique cluster pays mechanic for exciting gameplay"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Excellent graphics and authentic sounds"; New = "Stunning visual design with authentic Hawaiian theme"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Free Spins round with symbol drop mechanic"; New = "Free spins mechanism with symbol drop feature for better payouts"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Suitable for both casual players and high rollers"; New = "Suitable for both casual players and high rollers with wide betting range"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Not ideal for players who enjoy high-risk games with big payouts"; New = "Medium to low volatility may not appeal to high-risk players"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Smaller payouts due to lower volatility"; New = "Limited variety of symbols and gameplay features"; Style = $null; RPr = "" }
    [PSCustomObject]@{ Old = "Read our Aloha! Cluster Pays review and play for free. Learn strategies and enjoy excellent graphics, unique gameplay, and symbol drop mechanic in the free spins round."; New = "Experience the thrill of Aloha! Cluster Pays with free play and unique cluster pays mechanic."; Style = $null; RPr = "<w:rPr><w:i/></w:rPr>" }
)

# Build the concrete list of (start, end, newText, rPr) edits to perform.
$plan = @()
foreach ($edit in $edits) {
    foreach ($info in $paras) {
        if ($info.Text -ne $edit.Old) { continue }
        if (($edit.Style -ne $null) -and ($info.Style -ne $edit.Style)) { continue }
        $plan += [PSCustomObject]@{
            Start   = $info.Start
            End     = $info.End
            NewText = $edit.New
            RPr     = $edit.RPr
        }
    }
}

# Apply edits from the bottom of the document upward so that earlier
# (lower-offset) entries in $plan remain valid as later ones are applied.
$plan = $plan | Sort-Object -Property Start -Descending

foreach ($item in $plan) {
    Replace-ParaRun $d $item.Start $item.End $item.NewText $item.RPr
}

Write-Output "Applied $($plan.Count) replacements"
